$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = "'2"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'24.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = "'2"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.283"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = "'2"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.05795"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'2"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'6.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = "'2"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'3.146"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Value = "'2"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.8172"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Value = "'2"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.8563"
$ws.Range("D9").Style = "Normal"
$ws.Range("G9").Value = "'2"
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1361"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'9WazirXWRX"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'2"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.06951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'10MandalaExchangeTokenMDX"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'2"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03142"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'2"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.02873"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'12BitrueCoinBTR"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'2"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09399"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'13BitMartTokenBMX"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'2"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = "'MCDex"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'3.772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'14MCDexMCB"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'2"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.001519"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'15BitForexTokenBF"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'2"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "'CoinExToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.04682"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'16CoinExTokenCET"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'2"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = "'One"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.0005964"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'17OneONE"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'2"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.006288"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = "'2"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.001238"
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").Value = "'2"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.004608"
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").Value = "'2"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.00006906"
$ws.Range("D22").Style = "Normal"
$ws.Range("G22").Value = "'2"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'3.541"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Value = "'2"
$ws.Range("G23").Style = "Normal"
$ws.Range("G24").Value = "'2"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.3192"
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Value = "'2"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.1347"
$ws.Range("D26").Style = "Normal"
$ws.Range("G26").Value = "'2"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.1360"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = "'2"
$ws.Range("G27").Style = "Normal"
$ws.Range("D28").Value = "'0.0002331"
$ws.Range("D28").Style = "Normal"
$ws.Range("G28").Value = "'2"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'2"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'2"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'2"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'2"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'2"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'2"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'2"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'2"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'2"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'2"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'2"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.03655"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'2"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006265"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'2"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1054"
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").Value = "'2"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.002752"
$ws.Range("D43").Style = "Normal"
$ws.Range("G43").Value = "'2"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.007917"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = "'2"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005261"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = "'2"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("G46").Value = "'2"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.3702"
$ws.Range("D47").Style = "Normal"
$ws.Range("G47").Value = "'2"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.002373"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'47BOLOBOLOBestin24h"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'2"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = "'2"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("G50").Value = "'2"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'2"
$ws.Range("G51").Style = "Normal"
